# Logboek.xlsx - "pascal excel sheet updated"
#
# Adds Pascal's newest logboek entries (row 9-13), fills in the missing
# "Tijd (uur)" value for the "Microfoon werkend krijgen" activity (D6),
# and makes the Pascal tab the active/selected sheet (it was previously
# on Stijn's tab) with the cursor left on E9.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Pascal")

# Row 6 ("Microfoon werkend krijgen") gained a logged duration.
$ws.Range("D6").Value = 4

# New activities logged underneath the existing rows.
$ws.Range("B9").Value  = "nieuw protorype maken zonder mic"
$ws.Range("C9").Value  = 6
$ws.Range("D9").Value  = 1.25

$ws.Range("B10").Value = "MQtt onderzoek"
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 3

$ws.Range("B11").Value = "mqtt op esp32"
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 4

$ws.Range("B12").Value = "mqtt op anndoid"
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 3

$ws.Range("B13").Value = "esp32 bugs er uit halen"
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1

# Pascal's tab is now the one being worked in / left selected.
$ws.Select()
$ws.Range("E9").Select()
